# "commondataprovider for connecting excel"
# Adds a second worksheet ("validCredentialTest") with a valid-login test
# row (including a hyperlink to the expected landing URL), makes it the
# active/selected tab, and leaves the original "invalidCredentialTest"
# sheet's data untouched (only its tab-selected state changes because the
# new sheet becomes the active one).

$wb = $excel.ActiveWorkbook

# Insert the new sheet right after the existing "invalidCredentialTest"
# sheet so tab order is [invalidCredentialTest, validCredentialTest].
$firstSheet = $wb.Worksheets.Item(1)
$ws = $wb.Worksheets.Add($null, $firstSheet)
$ws.Name = "validCredentialTest"

# Header row
$ws.Range("A1").Value = "Username"
$ws.Range("B1").Value = "Password"
$ws.Range("C1").Value = "Expected Url"

# Data row - valid admin credentials + expected post-login URL
$ws.Range("A2").Value = "Admin"
$ws.Range("B2").Value = "admin123"
$ws.Range("C2").Value = "https://opensource-demo.orangehrmlive.com/index.php/dashboard"

# Turn the expected-URL cell into a real hyperlink (adds the Hyperlink
# cell style / font automatically, same as Excel does).
$ws.Hyperlinks.Add($ws.Range("C2"), "https://opensource-demo.orangehrmlive.com/index.php/dashboard") | Out-Null

# Widen column C so the URL is fully visible (best-fit-ish width).
$ws.Columns.Item(3).ColumnWidth = 56.8

# Match the authored selection/active-cell state on the new sheet.
$ws.Range("C3").Select() | Out-Null

# Make the new sheet the active tab.
$ws.Activate() | Out-Null
